$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 21415
$ws.Range("F3").Value = 3322
$ws.Range("F4").Value = 856
$ws.Range("G4").Value = 60
$ws.Range("F5").Value = 623
$ws.Range("F6").Value = 548
$ws.Range("F7").Value = 801
$ws.Range("F8").Value = 301
$ws.Range("F10").Value = 75
$ws.Range("F11").Value = 140
$ws.Range("F12").Value = 573
$ws.Range("F14").Value = 357
$ws.Range("F15").Value = 38
$ws.Range("F16").Value = 457
$ws.Range("F17").Value = 193
$ws.Range("F20").Value = 83
$ws.Range("F21").Value = 153

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 144
$ws.Range("F10").Value = 170

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6175
$ws.Range("F3").Value = 729
$ws.Range("F4").Value = 729
$ws.Range("F5").Value = 1735
$ws.Range("F6").Value = 87

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6176
$ws.Range("F3").Value = 729
$ws.Range("F4").Value = 729
$ws.Range("F5").Value = 1735
$ws.Range("F6").Value = 21415
$ws.Range("F7").Value = 3322
$ws.Range("F8").Value = 856
$ws.Range("G8").Value = 60
$ws.Range("F9").Value = 144
$ws.Range("F10").Value = 87
$ws.Range("F11").Value = 623
$ws.Range("F12").Value = 548
$ws.Range("F13").Value = 801
$ws.Range("F14").Value = 301
$ws.Range("F17").Value = 75
$ws.Range("F20").Value = 140
$ws.Range("F23").Value = 573
$ws.Range("F27").Value = 357
$ws.Range("F28").Value = 170
$ws.Range("F29").Value = 38
$ws.Range("F30").Value = 457
$ws.Range("F32").Value = 193
$ws.Range("F37").Value = 83
$ws.Range("F43").Value = 153

